$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44176
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 500
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("S2").Value = 1033
$ws.Range("T2").Value = 15
# Row 3
$ws.Range("D3").Value = 44565
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1111
# Row 4
$ws.Range("D4").Value = 44565
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 1000
# Row 5
$ws.Range("L5").Value = 'Segunda'
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("S5").Value = 889
# Row 6
$ws.Range("D6").Value = 44568
$ws.Range("L6").Value = 'Especial'
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 21000
$ws.Range("S6").Value = 1167
# Row 7
$ws.Range("D7").Value = 44568
$ws.Range("L7").Value = 'Primera'
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 1000
# Row 8
$ws.Range("D8").Value = 44568
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 889
$ws.Range("T8").Value = 18
# Row 9
$ws.Range("D9").Value = 44553
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 1111
$ws.Range("T9").Value = 18
# Row 10
$ws.Range("D10").Value = 44553
$ws.Range("M10").Value = 250
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("S10").Value = 1000
# Row 11
$ws.Range("D11").Value = 44553
$ws.Range("M11").Value = 250
# Row 12
$ws.Range("D12").Value = 44547
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 350
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 1111
# Row 13
$ws.Range("D13").Value = 44547
$ws.Range("M13").Value = 350
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 1000
# Row 14
$ws.Range("D14").Value = 44547
$ws.Range("M14").Value = 350
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 889
# Row 15
$ws.Range("D15").Value = 44551
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = '$/caja 18 kilos'
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 1111
$ws.Range("T15").Value = 18
# Row 16
$ws.Range("D16").Value = 44551
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 18
# Row 17
$ws.Range("D17").Value = 44551
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 889
# Row 18
$ws.Range("D18").Value = 44537
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 500
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 21000
$ws.Range("S18").Value = 1167
# Row 19
$ws.Range("D19").Value = 44537
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 250
$ws.Range("R19").Value = 'Región del Maule'
# Row 20
$ws.Range("D20").Value = 44159
$ws.Range("L20").Value = 'Tercera'
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 15500
$ws.Range("P20").Value = 15750
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 1050
# Row 24
$ws.Range("D24").Value = 44166
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 600
$ws.Range("N24").Value = 16000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 16500
$ws.Range("Q24").Value = '$/caja 15 kilos'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 1100
$ws.Range("T24").Value = 15
# Row 25
$ws.Range("D25").Value = 44530
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 500
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20500
$ws.Range("R25").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S25").Value = 1139
# Row 26
$ws.Range("D26").Value = 44187
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 350
$ws.Range("Q26").Value = '$/caja 15 kilos'
$ws.Range("S26").Value = 1067
$ws.Range("T26").Value = 15
# Row 27
$ws.Range("D27").Value = 44187
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 300
$ws.Range("N27").Value = 13000
$ws.Range("O27").Value = 13000
$ws.Range("P27").Value = 13000
$ws.Range("Q27").Value = '$/caja 15 kilos'
$ws.Range("S27").Value = 867
$ws.Range("T27").Value = 15
# Row 28
$ws.Range("D28").Value = 44540
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 600
$ws.Range("N28").Value = 16000
$ws.Range("O28").Value = 16000
$ws.Range("P28").Value = 16000
$ws.Range("R28").Value = 'Región del Maule'
$ws.Range("S28").Value = 889
# Row 29
$ws.Range("D29").Value = 44544
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 600
$ws.Range("N29").Value = 18000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19000
$ws.Range("S29").Value = 1056
# Row 30
$ws.Range("D30").Value = 44544
$ws.Range("L30").Value = 'Segunda'
$ws.Range("M30").Value = 300
$ws.Range("Q30").Value = '$/caja 18 kilos'
$ws.Range("S30").Value = 889
$ws.Range("T30").Value = 18
# Row 31
$ws.Range("D31").Value = 44169
$ws.Range("M31").Value = 500
$ws.Range("N31").Value = 15000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 15500
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1033
# Row 32
$ws.Range("D32").Value = 44194
$ws.Range("L32").Value = 'Segunda'
$ws.Range("M32").Value = 300
$ws.Range("N32").Value = 15000
$ws.Range("O32").Value = 16000
$ws.Range("P32").Value = 15500
$ws.Range("Q32").Value = '$/caja 15 kilos'
$ws.Range("S32").Value = 1033
$ws.Range("T32").Value = 15
# Row 33
$ws.Range("D33").Value = 44162
$ws.Range("L33").Value = 'Tercera'
$ws.Range("M33").Value = 500
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 16000
$ws.Range("P33").Value = 15500
$ws.Range("Q33").Value = '$/caja 15 kilos'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 1033
$ws.Range("T33").Value = 15
# Row 34
$ws.Range("D34").Value = 44533
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 350
$ws.Range("N34").Value = 24000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 24000
$ws.Range("R34").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S34").Value = 1333
# Row 35
$ws.Range("D35").Value = 44533
$ws.Range("M35").Value = 350
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 20000
$ws.Range("Q35").Value = '$/caja 18 kilos'
$ws.Range("R35").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S35").Value = 1111
$ws.Range("T35").Value = 18
# Row 36
$ws.Range("D36").Value = 44533
$ws.Range("L36").Value = 'Tercera'
$ws.Range("M36").Value = 350
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 17000
$ws.Range("P36").Value = 17000
$ws.Range("R36").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S36").Value = 944
